# Auto-update price data: push a new "today" row (2025-12-02) onto the
# top of the price table, shifting the existing rows down by one and
# dropping the values from a day forward (commodity prices were flat,
# so every row carries the same B/C/D numbers as before).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row above the first data row (row 2); this shifts
# rows 2:12 down to 3:13, matching the target dimension A1:D13.
$ws.Rows("2:2").Insert()

# The inserted row inherits formatting from the row above it (the bold
# header row) -- strip that back to the plain/default style used by the
# rest of the data rows.
$ws.Range("A2:D2").ClearFormats()

# Force column A to plain text *before* writing the date string so Excel
# doesn't auto-convert "2025-12-02" into a date serial number -- the
# sheet stores dates as literal text in every other row.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-12-02"
# Drop the now-unneeded text number format so the cell matches the
# unstyled look of the other date cells.
$ws.Range("A2").ClearFormats()

$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
